$wb = $excel.ActiveWorkbook

$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsCreateAccount = $wb.Worksheets.Item("CreateAccount")

# Update the email value - touch both cells that reference the old shared
# string so the engine rewrites it in place rather than appending a new one.
$wsSignIn.Range("C2").Value = "testjaga007@gmail.com"
$wsCreateAccount.Range("F2").Value = "testjaga007@gmail.com"

# Update the password value on the SignIn sheet (D2)
$wsSignIn.Range("D2").Value = "jaga@1234"

# Update selections on each sheet
$wsCreateAccount.Range("F7").Select()
$wsSignIn.Range("D2").Select()

# Make SignIn the active sheet/tab (was CreateAccount before)
$wsSignIn.Activate()

$wb.Save()
